# Actualización automática 2025-09-26 14:45:08
# Registra una venta de 142.56 (grupo "PIEDRA SINTERIZADA", septiembre) para
# el cliente CONZA VEGA FRANCO BLADYMIR del asesor RIOS CARRION ANGEL BENIGNO,
# y propaga el cambio a los totales / porcentajes dependientes.

$wb = $excel.ActiveWorkbook

# --- Hoja "VENTAS POR GRUPO" ---------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L9").Value = 142.56
$wsGrupo.Range("L26").Value = "2 de 24"

# --- Hoja "VENTA MENSUAL" -------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F9").Value = 142.56
$wsMensual.Range("F26").Value = 19699.59

# --- Hoja "CUMPLIMIENTO MENSUAL" -----------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D11").Value = 649.4400000000001
$wsCumplimiento.Range("E11").Value = 2272.78458185274
$wsCumplimiento.Range("F11").Value = 0.22224164563979

$wsCumplimiento.Range("D15").Value = 19699.59
$wsCumplimiento.Range("E15").Value = 38503.87623249458
$wsCumplimiento.Range("F15").Value = 0.3384607700391881
